$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 4) continuing from existing data in A2:B3
$ws.Range("A4").Value = -83.054317999999995
$ws.Range("B4").Value = 8.5982299999999992

# Update the selected cell to match the new active cell/selection (B5)
$ws.Range("B5").Select()
